# Replace pipe-delimited separators ("|") with comma-space separators (", ")
# in the "keywords" (D) and "junk-words" (E) columns of the report scraper
# example sheet, for easier parsing. One cell (D4) also drops a redundant
# trailing "fishing_report" entry, matching the upstream diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "fishing-reports, report"
    "D4"  = "montana-fly-fishing-reports, fishing-report"
    "D5"  = "fishing-reports, post"
    "E5"  = "news, tags, hashtags"
    "D6"  = "fishing-reports, fishing-report"
    "D7"  = "jackson-hole-fishing-report, fishing-report"
    "E7"  = "trip-report, tag"
    "D8"  = "montana-fishing-reports, montana-fishing-report"
    "D11" = "fishing-reports, fishing-report"
    "D13" = "reports, fishing-report"
    "D14" = "fishing-reports, fishing-report"
    "E14" = "uploads, page"
    "D15" = "fishing-reports-tips, fishing-reports, fishing-report"
    "E15" = "page, tag"
    "D19" = "fishing-report, fishing-reports"
    "D21" = "fishing-reports, fishing_report"
    "D24" = "fishing-report, river-report"
    "D25" = "montana-fishing-reports, fishing-report"
    "D26" = "river-reports, reports"
    "D27" = "madison-river-fishing-reports, fly-fishing-report"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
